# Insert a new weekly price record above row 23 ("Feria Lagunitas de
# Puerto Montt" / Espinaca sheet). This pushes the existing rows 23-41
# down to 24-42 (their contents are unchanged), growing the used range
# from A1:R41 to A1:R42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 23:41 down one row, leaving a blank row 23 to fill in.
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with the new record.
$ws.Cells.Item(23, 1).Value  = 4
$ws.Cells.Item(23, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(23, 3).Value  = "Los Lagos"
$ws.Cells.Item(23, 4).Value  = 44810
$ws.Cells.Item(23, 5).Value  = 10
$ws.Cells.Item(23, 6).Value  = 100112012
$ws.Cells.Item(23, 7).Value  = "Espinaca"
$ws.Cells.Item(23, 8).Value  = "Sin especificar"
$ws.Cells.Item(23, 9).Value  = "Primera"
$ws.Cells.Item(23, 10).Value = 30
$ws.Cells.Item(23, 11).Value = 12000
$ws.Cells.Item(23, 12).Value = 12000
$ws.Cells.Item(23, 13).Value = 12000
$ws.Cells.Item(23, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 16).Value = 1200
$ws.Cells.Item(23, 17).Value = 10
$ws.Cells.Item(23, 18).Value = "Hortaliza"
